# Auto-generated: refresh cryptocurrency price/volume data (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '97.334.35'
$ws.Range("E2").Value = '  +2.58%  '

# Row 3
$ws.Range("D3").Value = '3.593.99'
$ws.Range("E3").Value = '  +1.06%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.93'
$ws.Range("E5").Value = '  +2.27%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '655.57'
$ws.Range("E6").Value = '  +1.07%  '

# Row 7
$ws.Range("E7").Value = '  +16.29%  '

# Row 8
$ws.Range("E8").Value = '  +6.20%  '

# Row 9
$ws.Range("E9").Value = '  -0.07%  '

# Row 10
$ws.Range("E10").Value = '  +5.17%  '

# Row 11
$ws.Range("D11").Value = '3.591.73'
$ws.Range("E11").Value = '  +1.01%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.34'
$ws.Range("E12").Value = '  +5.00%  '

# Row 13
$ws.Range("E13").Value = '  +0.89%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.42'
$ws.Range("E14").Value = '  -0.74%  '

# Row 15
$ws.Range("D15").Value = '4.261.47'
$ws.Range("E15").Value = '  +1.11%  '

# Row 16
$ws.Range("D16").Value = '97.174.76'
$ws.Range("E16").Value = '  +2.37%  '

# Row 17
$ws.Range("E17").Value = '  +2.97%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.575.32'
$ws.Range("E18").Value = '  +0.67%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.81'
$ws.Range("E19").Value = '  +0.79%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.74'
$ws.Range("E20").Value = '  +1.73%  '

# Row 21
$ws.Range("E21").Value = '  +1.54%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.521'
$ws.Range("E22").Value = '  +9.39%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.51'
$ws.Range("E23").Value = '  +1.67%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '514.97'
$ws.Range("E24").Value = '  +2.00%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000205'
$ws.Range("E25").Value = '  +5.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.87'
$ws.Range("E26").Value = '  +1.52%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '101.19'
$ws.Range("E27").Value = '  +6.43%  '

# Row 28
$ws.Range("E28").Value = '  +4.87%  '

# Row 29
$ws.Range("D29").Value = '3.787.30'
$ws.Range("E29").Value = '  +1.07%  '

# Row 30
$ws.Range("E30").Value = '  +14.07%  '

# Row 31
$ws.Range("E31").Value = '  -0.07%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.84'
$ws.Range("E32").Value = '  +4.24%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.04%  '

# Row 34
$ws.Range("E34").Value = '  +3.64%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.72'
$ws.Range("E36").Value = '  +0.34%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '625.30'
$ws.Range("E37").Value = '  +7.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.81'
$ws.Range("E38").Value = '  +4.03%  '

# Row 39
$ws.Range("E39").Value = '  +1.91%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.65'
$ws.Range("E40").Value = '  +2.13%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.97'
$ws.Range("E41").Value = '  +11.76%  '

# Row 42
$ws.Range("E42").Value = '  +3.04%  '

# Row 43
$ws.Range("E43").Value = '  -0.03%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.01'
$ws.Range("E45").Value = '  +6.65%  '

# Row 46
$ws.Range("E46").Value = '  +7.55%  '

# Row 47
$ws.Range("E47").Value = '  +0.75%  '

# Row 48
$ws.Range("E48").Value = '  +1.13%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.407'
$ws.Range("E49").Value = '  +34.63%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.51'
$ws.Range("E50").Value = '  +5.82%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.06'
$ws.Range("E51").Value = '  -0.90%  '

